$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

$newNote = "Unsuccessful, reverted attempt to optimize context switch code"

# Row 19
$ws.Range("A19").Value = 41449
$ws.Range("A19").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B19").Value = 2
$ws.Range("D19").Value = $newNote

# Row 20
$ws.Range("A20").Value = 41450
$ws.Range("A20").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B20").Value = 1
$ws.Range("D20").Value = $newNote

$ws.Activate()
$ws.Range("A19").Select()
